$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 42249.75
$ws.Range("I47").Value = 18999
$ws.Range("K47").Value = 18999
$ws.Range("M47").Value = -18027

$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 5000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -4685
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 5000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -3908
$ws.Range("N79").ClearContents()

$ws.Range("H86").Value = 1778.5652
$ws.Range("I86").Value = 1550.3334
$ws.Range("J86").Value = 1925.2858
$ws.Range("K86").Value = 1550.3334
$ws.Range("L86").Value = 1925.2858
$ws.Range("M86").Value = -427.3334
$ws.Range("N86").Value = -4171.2858

$ws.Range("H89").Value = 1778.5652
$ws.Range("I89").Value = 1550.3334
$ws.Range("J89").Value = 1925.2858
$ws.Range("K89").Value = 7751.666999999999
$ws.Range("L89").Value = 9626.429
$ws.Range("M89").Value = -2135.666999999999
$ws.Range("N89").Value = -20858.429

$ws.Range("H107").Value = 1714.2174
$ws.Range("I107").Value = 1621.3572
$ws.Range("J107").Value = 1858.6666
$ws.Range("K107").Value = 1621.3572
$ws.Range("L107").Value = 1858.6666
$ws.Range("M107").Value = 298.6428000000001
$ws.Range("N107").Value = -5698.6666

$ws.Range("H118").Value = 935.7692
$ws.Range("J118").Value = 1011.5
$ws.Range("L118").Value = 3034.5
$ws.Range("N118").Value = -6348.5

$ws.Range("H129").Value = 911.25
$ws.Range("J129").Value = 919.6429000000001
$ws.Range("L129").Value = 2758.9287
$ws.Range("N129").Value = -12758.9287

$ws.Range("H137").Value = 3303.28
$ws.Range("I137").Value = 1519.5555
$ws.Range("K137").Value = 4558.666499999999
$ws.Range("M137").Value = -2008.666499999999

$ws.Range("H138").Value = 1746.1
$ws.Range("I138").Value = 591.2564
$ws.Range("J138").Value = 2484.4426
$ws.Range("K138").Value = 1773.7692
$ws.Range("L138").Value = 7453.327799999999
$ws.Range("M138").Value = 3366.2308
$ws.Range("N138").Value = -17733.3278

$ws.Range("H141").Value = 4374.7334
$ws.Range("I141").Value = 4524.885
$ws.Range("J141").Value = 3398.75
$ws.Range("K141").Value = 13574.655
$ws.Range("L141").Value = 10196.25
$ws.Range("M141").Value = -8394.655000000001
$ws.Range("N141").Value = -20556.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5501.8687
$ws.Range("I32").Value = 4751.143
$ws.Range("J32").Value = 8567.333000000001
$ws.Range("K32").Value = 4751.143
$ws.Range("L32").Value = 8567.333000000001
$ws.Range("M32").Value = -4464.143
$ws.Range("N32").Value = -9141.333000000001

$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H74").Value = 2757.7805
$ws.Range("I74").Value = 2562.342
$ws.Range("J74").Value = 5233.3335
$ws.Range("K74").Value = 2562.342
$ws.Range("L74").Value = 5233.3335
$ws.Range("M74").Value = -1688.342
$ws.Range("N74").Value = -6981.3335

$ws.Range("H77").Value = 2757.7805
$ws.Range("I77").Value = 2562.342
$ws.Range("J77").Value = 5233.3335
$ws.Range("K77").Value = 12811.71
$ws.Range("L77").Value = 26166.6675
$ws.Range("M77").Value = -8443.710000000001
$ws.Range("N77").Value = -34902.6675

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3309.2307
$ws.Range("I134").Value = 1840
$ws.Range("J134").Value = 5023.3335
$ws.Range("K134").Value = 5520
$ws.Range("L134").Value = 15070.0005
$ws.Range("M134").Value = -2985
$ws.Range("N134").Value = -20140.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 788.03705
$ws.Range("I113").Value = 664.5
$ws.Range("J113").Value = 886.86664
$ws.Range("K113").Value = 1993.5
$ws.Range("L113").Value = 2660.59992
$ws.Range("M113").Value = 176.5
$ws.Range("N113").Value = -7000.59992

$ws.Range("H122").Value = 3270.0667
$ws.Range("I122").Value = 1101
$ws.Range("J122").Value = 3603.7693
$ws.Range("K122").Value = 9909
$ws.Range("L122").Value = 32433.9237
$ws.Range("M122").Value = -7459
$ws.Range("N122").Value = -37333.9237

$ws.Range("H131").Value = 7576565.5
$ws.Range("I131").Value = 100000270
$ws.Range("J131").Value = 851.6393399999999
$ws.Range("K131").Value = 300000810
$ws.Range("L131").Value = 2554.91802
$ws.Range("M131").Value = -299995770
$ws.Range("N131").Value = -12634.91802

$ws.Range("H137").Value = 3343.4167
$ws.Range("I137").Value = 1372.5
$ws.Range("K137").Value = 4117.5
$ws.Range("M137").Value = 982.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 895.8
$ws.Range("I97").Value = 869.75
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 869.75
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -373.75
$ws.Range("N97").Value = -1992

$ws.Range("H132").Value = 2693.3125
$ws.Range("I132").Value = 1525.6842
$ws.Range("J132").Value = 4399.846
$ws.Range("K132").Value = 4577.0526
$ws.Range("L132").Value = 13199.538
$ws.Range("M132").Value = -2047.0526
$ws.Range("N132").Value = -18259.538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3163.2727
$ws.Range("I7").Value = 2363.3572
$ws.Range("J7").Value = 4563.125
$ws.Range("K7").Value = 2363.3572
$ws.Range("L7").Value = 4563.125
$ws.Range("M7").Value = -2251.3572
$ws.Range("N7").Value = -4787.125

$ws.Range("H61").Value = 1486.32
$ws.Range("I61").Value = 1412.25
$ws.Range("J61").Value = 1618
$ws.Range("K61").Value = 1412.25
$ws.Range("L61").Value = 1618
$ws.Range("M61").Value = -1210.25
$ws.Range("N61").Value = -2022

$ws.Range("H113").Value = 1486.32
$ws.Range("I113").Value = 1412.25
$ws.Range("J113").Value = 1618
$ws.Range("K113").Value = 1412.25
$ws.Range("L113").Value = 1618
$ws.Range("M113").Value = 757.75
$ws.Range("N113").Value = -5958

$ws.Range("H126").Value = 3163.2727
$ws.Range("I126").Value = 2363.3572
$ws.Range("J126").Value = 4563.125
$ws.Range("K126").Value = 7090.071599999999
$ws.Range("L126").Value = 13689.375
$ws.Range("M126").Value = -4620.071599999999
$ws.Range("N126").Value = -18629.375

$ws.Range("H132").Value = 3751.7407
$ws.Range("I132").Value = 1637.9445
$ws.Range("J132").Value = 7979.3335
$ws.Range("K132").Value = 4913.833500000001
$ws.Range("L132").Value = 23938.0005
$ws.Range("M132").Value = -2383.833500000001
$ws.Range("N132").Value = -28998.0005

$ws.Range("H136").Value = 2501.6538
$ws.Range("I136").Value = 908.41174
$ws.Range("J136").Value = 5511.1113
$ws.Range("K136").Value = 2725.23522
$ws.Range("L136").Value = 16533.3339
$ws.Range("M136").Value = -175.23522
$ws.Range("N136").Value = -21633.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 369.94736
$ws.Range("I113").Value = 300
$ws.Range("J113").Value = 420.81818
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 1262.45454
$ws.Range("M113").Value = 1270
$ws.Range("N113").Value = -5602.45454

$ws.Range("H126").Value = 2273.4546
$ws.Range("I126").Value = 1568.5
$ws.Range("K126").Value = 4705.5
$ws.Range("M126").Value = -2235.5

$ws.Range("H136").Value = 2650.4
$ws.Range("I136").Value = 810.5599999999999
$ws.Range("J136").Value = 7250
$ws.Range("K136").Value = 2431.68
$ws.Range("L136").Value = 21750
$ws.Range("M136").Value = 118.3200000000002
$ws.Range("N136").Value = -26850
